# Actualización automática de scrims_actualizado.xlsx (2025-07-23 22:40:26)
# Adds new scrim result rows to the "Open Business" and "Layer Cake" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Open Business": append rows 31-35 (dimension A3:N30 -> A3:N35)
# ---------------------------------------------------------------------
$wsOB = $wb.Worksheets.Item("Open Business")

# Row 31 needs the "Equipo 1" (style 5) look used by row 29; row 32 the
# same "Equipo 1" look; row 33 "Equipo 1" again; row 34 "Equipo 2" (style
# 4, like row 30); row 35 "Equipo 1" again. Copying a whole formatted row
# and pasting formats-only keeps the existing style indices instead of
# minting new ones.
$wsOB.Range("A29:N29").Copy()
$wsOB.Range("A31:N31").PasteSpecial(-4122)

$wsOB.Range("A29:N29").Copy()
$wsOB.Range("A32:N32").PasteSpecial(-4122)

$wsOB.Range("A29:N29").Copy()
$wsOB.Range("A33:N33").PasteSpecial(-4122)

$wsOB.Range("A30:N30").Copy()
$wsOB.Range("A34:N34").PasteSpecial(-4122)

$wsOB.Range("A29:N29").Copy()
$wsOB.Range("A35:N35").PasteSpecial(-4122)

# Row 31
$wsOB.Range("A31").Value = "GRAY"
$wsOB.Range("B31").Value = "MELODIE"
$wsOB.Range("C31").Value = "HANK"
$wsOB.Range("D31").Value = "DRACO"
$wsOB.Range("E31").Value = "KENJI"
$wsOB.Range("F31").Value = "LUMI"
$wsOB.Range("G31").Value = "Equipo 1"
$wsOB.Range("H31").Value = "TRB|Zeus 解開"
$wsOB.Range("I31").Value = "TRB|Lxffy"
$wsOB.Range("J31").Value = "TRB|R B M"
$wsOB.Range("K31").Value = "KCP|Fade"
$wsOB.Range("L31").Value = "KCP|Zoulan"
$wsOB.Range("M31").Value = "KCP|Tyrant"
$wsOB.Range("N31").Value = "20250723T203134.000Z"

# Row 32
$wsOB.Range("A32").Value = "GRAY"
$wsOB.Range("B32").Value = "MELODIE"
$wsOB.Range("C32").Value = "HANK"
$wsOB.Range("D32").Value = "DRACO"
$wsOB.Range("E32").Value = "KENJI"
$wsOB.Range("F32").Value = "LUMI"
$wsOB.Range("G32").Value = "Equipo 1"
$wsOB.Range("H32").Value = "TRB|Zeus 解開"
$wsOB.Range("I32").Value = "TRB|Lxffy"
$wsOB.Range("J32").Value = "TRB|R B M"
$wsOB.Range("K32").Value = "KCP|Fade"
$wsOB.Range("L32").Value = "KCP|Zoulan"
$wsOB.Range("M32").Value = "KCP|Tyrant"
$wsOB.Range("N32").Value = "20250723T202922.000Z"

# Row 33
$wsOB.Range("A33").Value = "MEEPLE"
$wsOB.Range("B33").Value = "ASH"
$wsOB.Range("C33").Value = "CORDELIUS"
$wsOB.Range("D33").Value = "HANK"
$wsOB.Range("E33").Value = "FINX"
$wsOB.Range("F33").Value = "JESSIE"
$wsOB.Range("G33").Value = "Equipo 1"
$wsOB.Range("H33").Value = "TRB|Zeus 解開"
$wsOB.Range("I33").Value = "TRB|R B M"
$wsOB.Range("J33").Value = "TRB|Lxffy"
$wsOB.Range("K33").Value = "KCP|Fade"
$wsOB.Range("L33").Value = "KCP|Tyrant"
$wsOB.Range("M33").Value = "KCP|Zoulan"
$wsOB.Range("N33").Value = "20250723T202204.000Z"

# Row 34
$wsOB.Range("A34").Value = "MEEPLE"
$wsOB.Range("B34").Value = "ASH"
$wsOB.Range("C34").Value = "CORDELIUS"
$wsOB.Range("D34").Value = "HANK"
$wsOB.Range("E34").Value = "FINX"
$wsOB.Range("F34").Value = "JESSIE"
$wsOB.Range("G34").Value = "Equipo 2"
$wsOB.Range("H34").Value = "TRB|Zeus 解開"
$wsOB.Range("I34").Value = "TRB|R B M"
$wsOB.Range("J34").Value = "TRB|Lxffy"
$wsOB.Range("K34").Value = "KCP|Fade"
$wsOB.Range("L34").Value = "KCP|Tyrant"
$wsOB.Range("M34").Value = "KCP|Zoulan"
$wsOB.Range("N34").Value = "20250723T202032.000Z"

# Row 35
$wsOB.Range("A35").Value = "MEEPLE"
$wsOB.Range("B35").Value = "ASH"
$wsOB.Range("C35").Value = "CORDELIUS"
$wsOB.Range("D35").Value = "HANK"
$wsOB.Range("E35").Value = "FINX"
$wsOB.Range("F35").Value = "JESSIE"
$wsOB.Range("G35").Value = "Equipo 1"
$wsOB.Range("H35").Value = "TRB|Zeus 解開"
$wsOB.Range("I35").Value = "TRB|R B M"
$wsOB.Range("J35").Value = "TRB|Lxffy"
$wsOB.Range("K35").Value = "KCP|Fade"
$wsOB.Range("L35").Value = "KCP|Tyrant"
$wsOB.Range("M35").Value = "KCP|Zoulan"
$wsOB.Range("N35").Value = "20250723T201744.000Z"

# ---------------------------------------------------------------------
# Sheet "Layer Cake": append row 43 (dimension A3:N42 -> A3:N43)
# ---------------------------------------------------------------------
$wsLC = $wb.Worksheets.Item("Layer Cake")

# Row 43 uses the "Equipo 2" (style 4) look, same as row 41.
$wsLC.Range("A41:N41").Copy()
$wsLC.Range("A43:N43").PasteSpecial(-4122)

$wsLC.Range("A43").Value = "ALLI"
$wsLC.Range("B43").Value = "MANDY"
$wsLC.Range("C43").Value = "CORDELIUS"
$wsLC.Range("D43").Value = "BUSTER"
$wsLC.Range("E43").Value = "CHESTER"
$wsLC.Range("F43").Value = "LOU"
$wsLC.Range("G43").Value = "Equipo 2"
$wsLC.Range("H43").Value = "TRB|R B M"
$wsLC.Range("I43").Value = "TRB|Zeus 解開"
$wsLC.Range("J43").Value = "TRB|Lxffy"
$wsLC.Range("K43").Value = "KCP|Fade"
$wsLC.Range("L43").Value = "KCP|Zoulan"
$wsLC.Range("M43").Value = "KCP|Tyrant"
$wsLC.Range("N43").Value = "20250723T203907.000Z"
